# [ADDITIONAL SCRAPING] added code to scrape more data about a player's
# batting performance in a match, also updated the excel sheets.
#
# 1. Insert a new "Player Info" sheet at the front of the workbook with
#    the player's basic info (ID, NAME, BATTING_HAND, BOWL_STYLE).
# 2. On "ODI Batting" and "ODI Bowling", rename the MATCH_CARD_LINK column
#    to MATCH_CODE and replace each scorecard URL with just the numeric
#    match code that was embedded in it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. New "Player Info" sheet, inserted before the existing first sheet.
# ---------------------------------------------------------------------
# NOTE: the sheet handle passed into Worksheets.Add(Before) ends up
# aliased to the *new* sheet once it is renamed, so re-fetch the
# "ODI Batting" sheet by name afterwards rather than reusing the old
# variable.
$infoSheet = $wb.Worksheets.Add($wb.Worksheets.Item("ODI Batting"))
$infoSheet.Name = "Player Info"

# Header row (bold / centered / bordered, same look as the other sheets).
$infoSheet.Cells.Item(1, 1).Value = "ID"
$infoSheet.Cells.Item(1, 2).Value = "NAME"
$infoSheet.Cells.Item(1, 3).Value = "BATTING_HAND"
$infoSheet.Cells.Item(1, 4).Value = "BOWL_STYLE"

$infoHeader = $infoSheet.Range("A1:D1")
$infoHeader.Font.Bold = $true
$infoHeader.HorizontalAlignment = -4108
$infoHeader.VerticalAlignment = -4160
$infoHeader.Borders.LineStyle = 1
$infoHeader.Borders.Weight = 2

# Data row. ID is numeric-looking text in the source data, so force the
# cell to stay text (matches how MATCH_NUMBER / RUNS_SCORED etc. are
# stored elsewhere in this workbook) before writing the value.
$idCell = $infoSheet.Cells.Item(2, 1)
$idCell.NumberFormat = "@"
$idCell.Value = "5663"

$infoSheet.Cells.Item(2, 2).Value = "Khushdil Shah"
$infoSheet.Cells.Item(2, 3).Value = "Left Handed"
$infoSheet.Cells.Item(2, 4).Value = "Does Not Bowl | Unknown"

# ---------------------------------------------------------------------
# 2. MATCH_CARD_LINK -> MATCH_CODE on "ODI Batting" (column D).
# ---------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Cells.Item(1, 4).Value = "MATCH_CODE"

$battingLastRow = $battingSheet.UsedRange.Rows.Count
for ($r = 2; $r -le $battingLastRow; $r++) {
    $cell = $battingSheet.Cells.Item($r, 4)
    $link = $cell.Value2
    if ($link) {
        $code = $link -replace '.*MatchCode=', ''
        $cell.NumberFormat = "@"
        $cell.Value = $code
    }
}

# ---------------------------------------------------------------------
# 3. MATCH_CARD_LINK -> MATCH_CODE on "ODI Bowling" (column B).
# ---------------------------------------------------------------------
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Cells.Item(1, 2).Value = "MATCH_CODE"

$bowlingLastRow = $bowlingSheet.UsedRange.Rows.Count
for ($r = 2; $r -le $bowlingLastRow; $r++) {
    $cell = $bowlingSheet.Cells.Item($r, 2)
    $link = $cell.Value2
    if ($link) {
        $code = $link -replace '.*MatchCode=', ''
        $cell.NumberFormat = "@"
        $cell.Value = $code
    }
}
